$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '82.198.92'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.89%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.198.24'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.69%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '625.39'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.70%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.293'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +21.73%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.587'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.191.82'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.594'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000261'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +12.39%  '
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('E14').Value = '  -3.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.784.71'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '32.10'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '82.011.81'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.96%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.188.36'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.65%  '
$ws.Range('E19').Value = '  +6.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.12'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '436.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.16'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.25'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.41'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +13.53%  '
$ws.Range('B26').Value = 'Aptos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.29'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.365.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '76.79'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.36%  '
$ws.Range('E30').Value = '  +3.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '591.32'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +11.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.14'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.53'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.147'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +22.20%  '
$ws.Range('E36').Value = '  +9.46%  '
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '22.91'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.18'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +11.05%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('E42').Value = '  +14.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.09'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +20.64%  '
$ws.Range('E44').Value = '  +3.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '161.15'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.26%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '188.57'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.54%  '
$ws.Range('E48').Value = '  +0.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '44.72'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '26.45'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.774'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.12%  '
